# Figure Controls readme: add four new rejection-toggle shortcuts and
# rename the old "Backspace" entry, moving the hidden _GoBack bookmark
# so it still marks the end of the user's last edit.

$d = $word.ActiveDocument

# --- Remove the existing _GoBack bookmark -----------------------------
# It currently sits inside the "Enter: Accept all trials (...)" paragraph,
# right before "but continue with original points)". It will be re-created
# later at its new location.
$oldMark = $d.Bookmarks.Item("_GoBack")
$oldMark.Delete()

# --- Insert the four new numbered rejection-toggle paragraphs ---------
# They go right before the "Backspace: Toggle flag for background
# movement detection" paragraph.
$backspacePara = $d.Paragraphs.Item(6)
$newParasBlock = "1: Toggle Background RMS rejection`r" + `
    "2: Toggle Background Voltage rejection`r" + `
    "3: Toggle MEP SD rejection`r" + `
    "4: Toggle MEP Voltage rejection`r"
$backspacePara.Range.InsertBefore($newParasBlock)

# --- Rewrite the old Backspace paragraph's text ------------------------
# A temporary marker is appended so the Find match below does not end
# exactly on the paragraph mark (placing a bookmark collapsed right at a
# paragraph boundary lands it in the wrong spot), then the marker is
# stripped once the bookmark has been anchored.
$d.Content.Find.Execute(
    "Backspace: Toggle flag for background movement detection",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Backspace: Clear all rejection criteria@@TMPMARK@@", 2)

# --- Re-create the _GoBack bookmark at the end of the new text ---------
$endRange = $d.Content
$endRange.Find.Execute("Backspace: Clear all rejection criteria")
$endRange.Collapse(0)
$d.Bookmarks.Add("_GoBack", $endRange)

# --- Strip the temporary marker -----------------------------------------
$d.Content.Find.Execute("@@TMPMARK@@", $true, $false, $false, $false, $false, $true, 1, $false, "", 2)
